# =====================================================================
# Adds 8 new "Error"/regression-diagnostic columns (R:Y) parsed out of the
# raw bracketed-array strings that used to just sit unused in column C,
# plus a scratch "Best Zscore" row (23/24) duplicating the best run (row 16)
# so it is easy to eyeball against the aggregate stats in rows 18-21.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers (R1:Y1) ------------------------------------------------
$ws.Range("R1").Value = "Y intercept"
$ws.Range("S1").Value = "Ereg"
$ws.Range("T1").Value = "Ep1"
$ws.Range("U1").Value = "Ep2"
$ws.Range("V1").Value = "Enp"
$ws.Range("W1").Value = "Epsi"
$ws.Range("X1").Value = "Ephi"
$ws.Range("Y1").Value = "Esa"

# --- Per-file data, rows 2-16 (columns R:Y) ------------------------------
# Row 2
$ws.Range("R2").Value = 1.3203125
$ws.Range("S2").Value = 2.69140625
$ws.Range("T2").Value = 3.453125
$ws.Range("U2").Value = 2.2060546900000002
$ws.Range("V2").Value = [double]"5.7617189999999999E-2"
$ws.Range("W2").Value = 0.13867188
$ws.Range("X2").Value = 0.67382812000000003
$ws.Range("Y2").Value = 0.36035156000000002

# Row 3 -- displayed in scientific notation, like the source data
$ws.Range("R3:Y3").NumberFormat = "0.00E+00"
$ws.Range("R3").Value = 1.28808594
$ws.Range("S3").Value = 3.6005859400000002
$ws.Range("T3").Value = 2.35546875
$ws.Range("U3").Value = 3.56054688
$ws.Range("V3").Value = 0.36230468799999999
$ws.Range("W3").Value = [double]"9.765625E-4"
$ws.Range("X3").Value = 0.83300781199999996
$ws.Range("Y3").Value = 0.32910156200000001

# Row 4
$ws.Range("R4").Value = 1.93164062
$ws.Range("S4").Value = 2.7607421900000002
$ws.Range("T4").Value = 2.296875
$ws.Range("U4").Value = [double]"9.375E-2"
$ws.Range("V4").Value = 0.79394531000000002
$ws.Range("W4").Value = [double]"9.7656199999999992E-3"
$ws.Range("X4").Value = 0.8828125
$ws.Range("Y4").Value = 0.26367188000000003

# Row 5
$ws.Range("R5").Value = 2.69726562
$ws.Range("S5").Value = 3.31054688
$ws.Range("T5").Value = 1.19726562
$ws.Range("U5").Value = 1.9140625
$ws.Range("V5").Value = 2.8681640599999998
$ws.Range("W5").Value = [double]"8.7890599999999996E-3"
$ws.Range("X5").Value = 1.03027344
$ws.Range("Y5").Value = 0.32226561999999997

# Row 6
$ws.Range("R6").Value = 2.5859375
$ws.Range("S6").Value = 2.82226562
$ws.Range("T6").Value = 2.1337890599999998
$ws.Range("U6").Value = 3.3798828099999998
$ws.Range("V6").Value = 1.85449219
$ws.Range("W6").Value = [double]"8.3007810000000001E-2"
$ws.Range("X6").Value = 0.58300781000000002
$ws.Range("Y6").Value = 0.28417968999999998

# Row 7
$ws.Range("R7").Value = 0.3125
$ws.Range("S7").Value = 3.28320312
$ws.Range("T7").Value = 0.6796875
$ws.Range("U7").Value = 3.98828125
$ws.Range("V7").Value = 0.38378906000000002
$ws.Range("W7").Value = [double]"6.8359379999999997E-2"
$ws.Range("X7").Value = 0.63476562000000003
$ws.Range("Y7").Value = 0.27636718999999998

# Row 8
$ws.Range("R8").Value = 1.71679688
$ws.Range("S8").Value = 3.74609375
$ws.Range("T8").Value = 2.6181640599999998
$ws.Range("U8").Value = 3.6640625
$ws.Range("V8").Value = [double]"9.765625E-2"
$ws.Range("W8").Value = [double]"2.4414060000000001E-2"
$ws.Range("X8").Value = 0.81054687999999997
$ws.Range("Y8").Value = 0.50488281000000002

# Row 9
$ws.Range("R9").Value = 3.51953125
$ws.Range("S9").Value = 3.84960938
$ws.Range("T9").Value = 0.79589843999999998
$ws.Range("U9").Value = 2.63476562
$ws.Range("V9").Value = 3.90625
$ws.Range("W9").Value = [double]"7.6171879999999997E-2"
$ws.Range("X9").Value = 0.96679687999999997
$ws.Range("Y9").Value = 0.35058593999999998

# Row 10
$ws.Range("R10").Value = 0.50195312000000003
$ws.Range("S10").Value = 2.57421875
$ws.Range("T10").Value = 2.96679688
$ws.Range("U10").Value = 2.92773438
$ws.Range("V10").Value = 1.21289062
$ws.Range("W10").Value = [double]"5.3710939999999999E-2"
$ws.Range("X10").Value = 0.8203125
$ws.Range("Y10").Value = 0.3046875

# Row 11 -- displayed in scientific notation, like the source data
$ws.Range("R11:Y11").NumberFormat = "0.00E+00"
$ws.Range("R11").Value = 2.484375
$ws.Range("S11").Value = 3.3876953099999998
$ws.Range("T11").Value = 2.5205078099999998
$ws.Range("U11").Value = 3.8125
$ws.Range("V11").Value = 2.34375
$ws.Range("W11").Value = [double]"2.9296875E-3"
$ws.Range("X11").Value = 1.17773438
$ws.Range("Y11").Value = 0.404296875

# Row 12
$ws.Range("R12").Value = 3.2421875
$ws.Range("S12").Value = 3.73242188
$ws.Range("T12").Value = 3.56640625
$ws.Range("U12").Value = 3.76953125
$ws.Range("V12").Value = 0.69824218999999998
$ws.Range("W12").Value = [double]"5.5664060000000001E-2"
$ws.Range("X12").Value = 1.02929688
$ws.Range("Y12").Value = 0.37304688000000003

# Row 13
$ws.Range("R13").Value = 3.515625
$ws.Range("S13").Value = 3.3212890599999998
$ws.Range("T13").Value = 1.92578125
$ws.Range("U13").Value = 3.4697265599999998
$ws.Range("V13").Value = 0.87402343999999998
$ws.Range("W13").Value = [double]"1.7578119999999999E-2"
$ws.Range("X13").Value = 0.76660156000000002
$ws.Range("Y13").Value = 0.44824218999999998

# Row 14
$ws.Range("R14").Value = 3.96484375
$ws.Range("S14").Value = 2.7470703099999998
$ws.Range("T14").Value = 0.12890625
$ws.Range("U14").Value = 2.9365234400000002
$ws.Range("V14").Value = 0.25878906000000002
$ws.Range("W14").Value = [double]"8.3984379999999997E-2"
$ws.Range("X14").Value = 0.8125
$ws.Range("Y14").Value = 0.2578125

# Row 15
$ws.Range("R15").Value = 2.7177734400000002
$ws.Range("S15").Value = 3.9873046900000002
$ws.Range("T15").Value = 1.35546875
$ws.Range("U15").Value = 3.9638671900000002
$ws.Range("V15").Value = 3.0009765599999998
$ws.Range("W15").Value = [double]"5.3710939999999999E-2"
$ws.Range("X15").Value = 0.83007812000000003
$ws.Range("Y15").Value = 0.42871093999999998

# Row 16
$ws.Range("R16").Value = 3.5185546900000002
$ws.Range("S16").Value = 3.1015625
$ws.Range("T16").Value = 0.22753905999999999
$ws.Range("U16").Value = 1.4609375
$ws.Range("V16").Value = 1.953125
$ws.Range("W16").Value = 0.12890625
$ws.Range("X16").Value = 1.05175781
$ws.Range("Y16").Value = 0.29785156000000002

# --- Aggregate formulas for the new columns, mirroring D:Q ---------------
$ws.Range("R18:Y18").Formula = "=AVERAGE(R2:R16)"
$ws.Range("R19:Y19").Formula = "=STDEV(R2:R16)"
$ws.Range("R20:Y20").Formula = "=MAX(R2:R16)"
$ws.Range("R21:Y21").Formula = "=MIN(R2:R16)"

# --- Scratch "Best Zscore" block: row 16 (best run) duplicated for easy --
# --- side-by-side comparison against the Min/Max/Mean rows above ---------
$ws.Range("A23").Value = "Best Zscore"
$ws.Range("A24").Value = $ws.Range("A16").Value2
$ws.Range("C24").Value = $ws.Range("C16").Value2
$ws.Range("B24").Value = $ws.Range("B16").Value2
$ws.Range("D24").Value = $ws.Range("D16").Value2
$ws.Range("E24").Value = $ws.Range("E16").Value2
$ws.Range("F24").Value = $ws.Range("F16").Value2
$ws.Range("G24").Value = $ws.Range("G16").Value2
$ws.Range("H24").Value = $ws.Range("H16").Value2
$ws.Range("I24").Value = $ws.Range("I16").Value2
$ws.Range("J24").Value = $ws.Range("J16").Value2
$ws.Range("K24").Value = $ws.Range("K16").Value2
$ws.Range("L24").Value = $ws.Range("L16").Value2
$ws.Range("M24").Value = $ws.Range("M16").Value2
$ws.Range("N24").Value = $ws.Range("N16").Value2
$ws.Range("O24").Value = $ws.Range("O16").Value2
$ws.Range("P24").Value = $ws.Range("P16").Value2
$ws.Range("Q24").Value = $ws.Range("Q16").Value2
$ws.Range("R24").Value = $ws.Range("R16").Value2
$ws.Range("S24").Value = $ws.Range("S16").Value2
$ws.Range("T24").Value = $ws.Range("T16").Value2
$ws.Range("U24").Value = $ws.Range("U16").Value2
$ws.Range("V24").Value = $ws.Range("V16").Value2
$ws.Range("W24").Value = $ws.Range("W16").Value2
$ws.Range("X24").Value = $ws.Range("X16").Value2
$ws.Range("Y24").Value = $ws.Range("Y16").Value2

# --- Column A is now a hair too narrow for the longer file names; resize -
$ws.Columns("A:A").ColumnWidth = 15.25

# --- Selection left on the header row after the edit ---------------------
$ws.Rows("1:1").Select()
